$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('C2').Value = 'Unknown Title'
$ws.Range('D2').Value = 'Unknown Abstract'
$ws.Range('E2').Value = '[]'
$ws.Range('F2').Value = 'not found'
$ws.Range('G2').Value = 'N/A'
$ws.Range('H2').Value = '1970-01-01'
$ws.Range('J2').Value = ''
$ws.Range('C3').Value = 'Unknown Title'
$ws.Range('D3').Value = 'Unknown Abstract'
$ws.Range('E3').Value = '[]'
$ws.Range('F3').Value = 'not found'
$ws.Range('G3').Value = 'N/A'
$ws.Range('H3').Value = '1970-01-01'
$ws.Range('I3').Value = ''
$ws.Range('C4').Value = 'Unknown Title'
$ws.Range('D4').Value = 'Unknown Abstract'
$ws.Range('E4').Value = '[]'
$ws.Range('F4').Value = 'not found'
$ws.Range('G4').Value = 'N/A'
$ws.Range('H4').Value = '1970-01-01'
$ws.Range('I4').Value = ''
$ws.Range('E5').Value = '[Tuo%Ji%NULL%0, Hai-Lian%Chen%NULL%1, Jing%Xu%NULL%0, Ling-Ning%Wu%NULL%1, Jie-Jia%Li%NULL%1, Kai%Chen%NULL%1, Gang%Qin%tonygqin@ntu.edu.cn%1]'
$ws.Range('I5').Value = ''
$ws.Range('J5').Value = 'Oxford University Press'
$ws.Range('C6').Value = 'Unknown Title'
$ws.Range('E6').Value = '[]'
$ws.Range('F6').Value = 'not found'
$ws.Range('G6').Value = 'N/A'
$ws.Range('H6').Value = '1970-01-01'
$ws.Range('J6').Value = ''
$ws.Range('E7').Value = '[Hien%Lau%hlau2@uci.edu%0, Veria%Khosrawipour%veriakhosrawipour@yahoo.de%1, Piotr%Kocbach%piotr.kocbach@uwm.edu.pl%1, Agata%Mikolajczyk%agata.mikolajczyk@upwr.edu.pl%1, Justyna%Schubert%justyna.schubert@upwr.edu.pl%1, Jacek%Bania%jacek.bania@upwr.edu.pl%1, Tanja%Khosrawipour%tkhosrawipour@gmail.com%1]'
$ws.Range('I7').Value = ''
$ws.Range('J7').Value = 'Oxford University Press'
$ws.Range('D8').Value = 'Since the beginning of the COVID-19 epidemic in Italy, the Italian Government implemented several restrictive measures to contain the spread of the infection.
 Data shows that, among these measures, the lockdown implemented as of 9 March had a positive impact, in particular the central and southern regions of Italy, while other actions appeared to be less effective.
 When the true prevalence of a disease is unknown, it is possible estimate it, based on mortality data and the assumptive case-fatality rate of the disease.
 Given these assumptions, the estimated period-prevalence of COVID-19 in Italy varies from 0.35% in Sicily to 13.3% in Lombardy.
'
$ws.Range('E8').Value = '[Carlo%Signorelli%NULL%0, Thea%Scognamiglio%NULL%1, Anna%Odone%NULL%0]'
$ws.Range('I8').Value = ''
$ws.Range('J8').Value = 'Mattioli 1885'
$ws.Range('C9').Value = 'Unknown Title'
$ws.Range('D9').Value = 'Unknown Abstract'
$ws.Range('E9').Value = '[]'
$ws.Range('F9').Value = 'not found'
$ws.Range('G9').Value = 'N/A'
$ws.Range('H9').Value = '1970-01-01'
$ws.Range('I9').Value = ''
$ws.Range('D10').Value = 'From the end of February, the SARS-CoV-2 epidemic in Spain has been following the footsteps of that in Italy very closely.
 We have analyzed the trends of incident cases, deaths, and intensive care unit admissions (ICU) in both countries before and after their respective national lockdowns using an interrupted time-series design.
 Data was analyzed with quasi-Poisson regression using an interaction model to estimate the change in trends.
 After the first lockdown, incidence trends were considerably reduced in both countries.
 However, although the slopes have been flattened for all outcomes, the trends kept rising.
 During the second lockdown, implementing more restrictive measures for mobility, it has been a change in the trend slopes for both countries in daily incident cases and ICUs.
 This improvement indicates that the efforts overtaken are being successful in flattening the epidemic curve, and reinforcing the belief that we must hold on.
'
$ws.Range('E10').Value = '[Aurelio%Tobías%NULL%0]'
$ws.Range('I10').Value = ''
$ws.Range('J10').Value = 'Published by Elsevier B.V.'
$ws.Range('D11').Value = 'Coronavirus disease 2019 (COVID-19) is a global pandemic.
 Non-pharmacological interventions, such as lockdown and mass testing, remain as the mainstay of control measures for the outbreak.
 We aim to evaluate the effectiveness of mass testing, lockdown, or a combination of both to control COVID-19 pandemic.
 A systematic search on 11 major databases was conducted on June 8, 2020. This review is registered in Prospero (CRD420201 90546).
 We included primary studies written in English which investigate mass screening, lockdown, or a combination of both to control and/or mitigate the COVID-19 pandemic.
 There are four important outcomes as selected by WHO experts for their decision- making process: incident cases, onward transmission, mortality, and resource use.
 Among 623 studies, only 14 studies met our criteria.
 Four observational studies were rated as strong evidence and ten modelling studies were rated as moderate evidence.
 Based on one modelling study, mass testing reduced the total infected people compared to no mass testing.
 For lockdown, ten studies consistently showed that it successfully reduced the incidence, onward transmission, and mortality rate of COVID-19. A limited evidence showed that a combination of lockdown and mass screening resulted in a greater reduction of incidence and mortality rate compared to lockdown only.
 However, there is not enough evidence on the effectiveness of mass testing only.
'
$ws.Range('E11').Value = '[Nadya%Johanna%NULL%0, Henrico%Citrawijaya%NULL%1, Grace%Wangge%NULL%1]'
$ws.Range('I11').Value = ''
$ws.Range('J11').Value = 'PAGEPress Publications, Pavia, Italy'
$ws.Range('C12').Value = 'Unknown Title'
$ws.Range('E12').Value = '[]'
$ws.Range('F12').Value = 'not found'
$ws.Range('G12').Value = 'N/A'
$ws.Range('H12').Value = '1970-01-01'
$ws.Range('J12').Value = ''
$ws.Range('D13').Value = '
              •
              Since January 23rd 2020, stringent measures for controlling the novel coronavirus epidemics have been enforced and strengthened in mainland China.
'
$ws.Range('E13').Value = '[Biao%Tang%NULL%0, Fan%Xia%NULL%1, Sanyi%Tang%NULL%1, Nicola Luigi%Bragazzi%NULL%1, Qian%Li%NULL%0, Xiaodan%Sun%NULL%1, Juhua%Liang%NULL%1, Yanni%Xiao%yxiao@mail.xjtu.edu.cn%1, Jianhong%Wu%wujh@yorku.ca%1]'
$ws.Range('I13').Value = ''
$ws.Range('J13').Value = 'The Author(s). Published by Elsevier Ltd on behalf of International Society for Infectious Diseases.'
